$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 1: "Medicao Final" value changes from 107 -> 104
# -----------------------------------------------------------------
$t2 = $d.Tables.Item(2)
$cell1 = $t2.Cell(3, 2)
$cell1.Range.Find.Execute("107", $false, $false, $false, $false, $false, $true, 1, $false, "104", 2)

# Word drops its "_GoBack" last-edit bookmark right where the user
# finished typing. That spot is also the very last character of the
# paragraph, which this host mishandles when a *collapsed* range is
# built right on a paragraph-end boundary -- so nudge the boundary
# out of the way first with a throw-away character, plant the
# bookmark, then remove the scratch character again.
$p1 = $cell1.Range.Paragraphs.Item(1)
$editEnd = $p1.Range.End - 1
$d.Range($editEnd, $editEnd).InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $d.Range($editEnd, $editEnd))
$d.Range($editEnd, $editEnd + 1).Delete()

# -----------------------------------------------------------------
# Edit 2: the directory reference is retyped as one continuous run,
# "...01_OS4776\02_TESTES.", so the old "_GoBack" bookmark that used
# to sit mid-word there goes away.
# -----------------------------------------------------------------
$t5 = $d.Tables.Item(5)
$cell2 = $t5.Cell(18, 2)
$cell2.Range.Find.Execute("\02_TESTES.", $false, $false, $false, $false, $false, $true, 1, $false, "_TEMP_RETYPE_", 2)
$cell2.Range.Find.Execute("_TEMP_RETYPE_", $false, $false, $false, $false, $false, $true, 1, $false, "\02_TESTES.", 2)
